$d = $word.ActiveDocument

$pairs = @(
    @("177÷6=29, 3", "676÷7=96, 4"),
    @("753÷9=83, 6", "863÷4=215, 3"),
    @("239÷8=29, 7", "563÷9=62, 5"),
    @("616÷7=88, 0", "146÷4=36, 2"),
    @("806÷9=89, 5", "756÷8=94, 4"),
    @("925÷6=154, 1", "520÷4=130, 0"),
    @("810÷7=115, 5", "868÷3=289, 1"),
    @("168÷2=84, 0", "551÷5=110, 1"),
    @("185÷8=23, 1", "531÷5=106, 1"),
    @("448÷7=64, 0", "839÷8=104, 7"),
    @("433÷4=108, 1", "487÷8=60, 7"),
    @("748÷9=83, 1", "587÷6=97, 5"),
    @("313÷4=78, 1", "584÷8=73, 0"),
    @("815÷3=271, 2", "943÷5=188, 3"),
    @("801÷9=89, 0", "980÷8=122, 4"),
    @("194÷6=32, 2", "664÷4=166, 0"),
    @("918÷5=183, 3", "879÷8=109, 7"),
    @("759÷6=126, 3", "718÷3=239, 1"),
    @("483÷9=53, 6", "773÷5=154, 3"),
    @("808÷8=101, 0", "669÷7=95, 4"),
    @("691÷3=230, 1", "272÷7=38, 6"),
    @("946÷3=315, 1", "684÷4=171, 0"),
    @("140÷8=17, 4", "312÷5=62, 2"),
    @("194÷4=48, 2", "168÷4=42, 0"),
    @("726÷6=121, 0", "531÷5=106, 1")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
